# Updated Excel file with timestamp via Streamlit
#
# The sheet originally held 3 societies (Soleo Health, University of Miami,
# Dava Oncology LP) across rows 2-4. The new run only reports on
# "Soleo Health" (row 2) with refreshed answers, drops the other two rows,
# and appends a new "Last Updated" timestamp column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Write "1500" into B2 as TEXT (not a number) --------------------------
# Assigning a numeric-looking string straight to .Value auto-converts it to
# a real number. Stage it in a scratch cell that is pre-formatted as Text,
# then copy/paste-values it into B2 so B2 keeps the default (unstyled) cell
# format while still holding a text value. The scratch cell lives in row 3,
# which gets deleted later, so no trace of it survives in the final sheet.
$scratch = $ws.Range("Z3")
$scratch.NumberFormat = "@"
$scratch.Value = "1500"
$scratch.Copy()
$ws.Range("B2").PasteSpecial(-4163) | Out-Null

# --- Refresh the rest of row 2 with the new Soleo Health answers ----------
$ws.Range("C2").Value = "No, Soleo Health does not encompass community sites. Soleo Health is a specialty pharmacy and infusion services provider focused on patient care in the home setting, rather than community sites."
$ws.Range("D2").Value = "No, Soleo Health is not influential on state or local policy. Soleo Health operates as a healthcare provider, focusing on specialty infusion services rather than advocacy or policy-making activities."
$ws.Range("E2").Value = "No, Soleo Health, The organization does not provide engagement opportunities with leadership. The company may lack transparency or structured programs to facilitate interactions with its leadership team."
$ws.Range("F2").Value = "No, Soleo Health does not provide support for clinical trial recruitment. Soleo Health focuses on specialty infusion services for patients."
$ws.Range("G2").Value = "No, Soleo Health does not provide engaging opportunities with payors. Payor engagement is not a prominent focus for Soleo Health."
$ws.Range("H2").Value = "No, Soleo Health does not include area experts on its board. The company primarily focuses on providing specialty infusion services and healthcare solutions rather than having industry-specific experts on its board of directors."
$ws.Range("I2").Value = "no, There is no public information available on therapeutic research collaborations involving Soleo Health."
$ws.Range("J2").Value = "No, justification: There is no public information available on the composition of Soleo Health's board to confirm if it includes top therapeutic area experts."
$ws.Range("K2").Value = "Texas"

# --- Drop the University of Miami (row 3) and Dava Oncology (row 4) rows --
$ws.Rows("3:4").Delete()

# --- Add the new "Last Updated" column with the run timestamp -------------
$ws.Range("L1").Value = "Last Updated"
# Match the bold/centered/bordered header formatting used by the other
# header cells (copy format only from K1, an existing header cell).
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122) | Out-Null
$ws.Range("L1").Value = "Last Updated"

$ws.Range("L2").Value = "2025-03-12 12:46:35"

$excel.CutCopyMode = 0
